$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 21.90961833333333
$ws.Range("H2").Value = 65.728855
$ws.Range("I2").Value = 0.6422049224355262
$ws.Range("J2").Value = 0.6422049224355262
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 13.51863466666667
$ws.Range("N2").Value = 40.555904
$ws.Range("O2").Value = 0.5370474071126665
$ws.Range("P2").Value = 0.5370474071126665
$ws.Range("Q2").Value = 296.1881259344355
$ws.Range("R2").Value = 2665.69313340992
$ws.Range("S2").Value = 0.3448944884289905
$ws.Range("T2").Value = 0.3448944884289905
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 21.90961833333333
$ws.Range("H3").Value = 65.728855
$ws.Range("I3").Value = 0.6422049224355262
$ws.Range("J3").Value = 0.6422049224355262
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 11.419192
$ws.Range("N3").Value = 34.257576
$ws.Range("O3").Value = 0.4536439963159277
$ws.Range("P3").Value = 0.4536439963159277
$ws.Range("Q3").Value = 250.1901383950533
$ws.Range("R3").Value = 2251.71124555548
$ws.Range("S3").Value = 0.2913324074674125
$ws.Range("T3").Value = 0.2913324074674125
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 21.90961833333333
$ws.Range("H4").Value = 65.728855
$ws.Range("I4").Value = 0.6422049224355262
$ws.Range("J4").Value = 0.6422049224355262
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.2343173333333333
$ws.Range("N4").Value = 0.702952
$ws.Range("O4").Value = 0.00930859657140581
$ws.Range("P4").Value = 0.00930859657140581
$ws.Range("Q4").Value = 5.133803342217777
$ws.Range("R4").Value = 46.20423007996
$ws.Range("S4").Value = 0.005978026539123274
$ws.Range("T4").Value = 0.005978026539123274
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 10.79980633333333
$ws.Range("H5").Value = 32.399419
$ws.Range("I5").Value = 0.3165590875704608
$ws.Range("J5").Value = 0.3165590875704607
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 13.51863466666667
$ws.Range("N5").Value = 40.555904
$ws.Range("O5").Value = 0.5370474071126665
$ws.Range("P5").Value = 0.5370474071126665
$ws.Range("Q5").Value = 145.9986362910862
$ws.Range("R5").Value = 1313.987726619776
$ws.Range("S5").Value = 0.1700072371776675
$ws.Range("T5").Value = 0.1700072371776675
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 10.79980633333333
$ws.Range("H6").Value = 32.399419
$ws.Range("I6").Value = 0.3165590875704608
$ws.Range("J6").Value = 0.3165590875704607
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 11.419192
$ws.Range("N6").Value = 34.257576
$ws.Range("O6").Value = 0.4536439963159277
$ws.Range("P6").Value = 0.4536439963159277
$ws.Range("Q6").Value = 123.3250620831493
$ws.Range("R6").Value = 1109.925558748344
$ws.Range("S6").Value = 0.1436051295555875
$ws.Range("T6").Value = 0.1436051295555875
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 10.79980633333333
$ws.Range("H7").Value = 32.399419
$ws.Range("I7").Value = 0.3165590875704608
$ws.Range("J7").Value = 0.3165590875704607
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.2343173333333333
$ws.Range("N7").Value = 0.702952
$ws.Range("O7").Value = 0.00930859657140581
$ws.Range("P7").Value = 0.00930859657140581
$ws.Range("Q7").Value = 2.530581820543111
$ws.Range("R7").Value = 22.775236384888
$ws.Range("S7").Value = 0.002946720837205743
$ws.Range("T7").Value = 0.002946720837205742
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.406817
$ws.Range("H8").Value = 4.220451
$ws.Range("I8").Value = 0.04123598999401312
$ws.Range("J8").Value = 0.04123598999401312
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 13.51863466666667
$ws.Range("N8").Value = 40.555904
$ws.Range("O8").Value = 0.5370474071126665
$ws.Range("P8").Value = 0.5370474071126665
$ws.Range("Q8").Value = 19.018245065856
$ws.Range("R8").Value = 171.164205592704
$ws.Range("S8").Value = 0.02214568150600861
$ws.Range("T8").Value = 0.0221456815060086
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.406817
$ws.Range("H9").Value = 4.220451
$ws.Range("I9").Value = 0.04123598999401312
$ws.Range("J9").Value = 0.04123598999401312
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 11.419192
$ws.Range("N9").Value = 34.257576
$ws.Range("O9").Value = 0.4536439963159277
$ws.Range("P9").Value = 0.4536439963159277
$ws.Range("Q9").Value = 16.064713431864
$ws.Range("R9").Value = 144.582420886776
$ws.Range("S9").Value = 0.01870645929292772
$ws.Range("T9").Value = 0.01870645929292772
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.406817
$ws.Range("H10").Value = 4.220451
$ws.Range("I10").Value = 0.04123598999401312
$ws.Range("J10").Value = 0.04123598999401312
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.2343173333333333
$ws.Range("N10").Value = 0.702952
$ws.Range("O10").Value = 0.00930859657140581
$ws.Range("P10").Value = 0.00930859657140581
$ws.Range("Q10").Value = 0.329641607928
$ws.Range("R10").Value = 2.966774471352
$ws.Range("S10").Value = 0.0003838491950767949
$ws.Range("T10").Value = 0.0003838491950767948
